$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "(주)코그넷나인"
$ws.Range("B9").Value = "프론트엔드개발자"
$ws.Range("C9").Value = "https://www.jobplanet.co.kr/job/search?posting_ids%5B%5D=1290844"
$ws.Range("D9").Value = "경력"
$ws.Range("E9").Value = "react"

$ws.Range("A10").Value = "(주)앤씨앤"
$ws.Range("B10").Value = "Embedded Application 개발자 채용"
$ws.Range("C10").Value = "https://www.jobplanet.co.kr/job/search?posting_ids%5B%5D=1290660"
$ws.Range("D10").Value = "경력"

$ws.Range("A11").Value = "(주)유진로봇"
$ws.Range("B11").Value = "개발본부 AMS팀 백엔드(Back-end) 개발자"
$ws.Range("C11").Value = "https://www.jobplanet.co.kr/job/search?posting_ids%5B%5D=1269944"
$ws.Range("D11").Value = "경력"
$ws.Range("E11").Value = "javascript"
